$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column (05-dec) before the old EH column ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Insert a new blank column at EH, shifting the "01-oct." ... "31-oct." block one column right
$wsPrix.Range("EH:EH").Insert()

# New column header
$wsPrix.Range("EH1").Value = "05-dec"

# New column data for rows 2..25 is all "-" (no data available, like the other forward columns)
$wsPrix.Range("EH2:EH25").Value = "-"

# --- Sheet "Gaz": append new row 168 ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A168").NumberFormat = "@"
$wsGaz.Range("A168").Value = "2025-12-03"
$wsGaz.Range("A168").Style = "Normal"
$wsGaz.Range("B168").Value = 27.05

# --- Sheet "CO2": append new row 168 ---
$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Range("A168").NumberFormat = "@"
$wsCO2.Range("A168").Value = "2025-12-03"
$wsCO2.Range("A168").Style = "Normal"
$wsCO2.Range("B168").Value = 81.35
